$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remember the width of column B (household_id) so the newly inserted
# column can reuse it - matches the target layout where the new
# collector_id column has the same width as its neighbours.
$hhIdWidth = $ws.Columns("B").ColumnWidth

# Insert a new column before column C (collector_id), shifting existing
# household_size..national_id columns one to the right.
$ws.Columns("C").Insert()
$ws.Columns("C").ColumnWidth = $hhIdWidth

# New header for the inserted column
$ws.Range("C1").Value = "collector_id"

# New data values for the inserted column
$ws.Range("C2").Value = "IND-24-0000.0012"
$ws.Range("C3").Value = "IND-24-0000.0013"
